# Update "想去人数" (want-to-go count) figures in column F across all sheets
# to match the latest scrape (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 400
$ws.Range("F4").Value = 1356
$ws.Range("F5").Value = 248
$ws.Range("F6").Value = 2613
$ws.Range("F7").Value = 984
$ws.Range("F8").Value = 19055
$ws.Range("F9").Value = 66
$ws.Range("F10").Value = 2082
$ws.Range("F11").Value = 701
$ws.Range("F12").Value = 609
$ws.Range("F13").Value = 374
$ws.Range("F14").Value = 637
$ws.Range("F15").Value = 211
$ws.Range("F16").Value = 222
$ws.Range("F18").Value = 334
$ws.Range("F19").Value = 55
$ws.Range("F20").Value = 226
$ws.Range("F22").Value = 141
$ws.Range("F25").Value = 84

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 5
$ws.Range("F7").Value = 6
$ws.Range("F9").Value = 247
$ws.Range("F10").Value = 247
$ws.Range("F16").Value = 76
$ws.Range("F18").Value = 40

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5940
$ws.Range("F3").Value = 606

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5940
$ws.Range("F4").Value = 606
$ws.Range("F6").Value = 400
$ws.Range("F8").Value = 1356
$ws.Range("F10").Value = 248
$ws.Range("F13").Value = 2613
$ws.Range("F14").Value = 984
$ws.Range("F15").Value = 19055
$ws.Range("F16").Value = 5
$ws.Range("F17").Value = 6
$ws.Range("F18").Value = 66
$ws.Range("F20").Value = 247
$ws.Range("F21").Value = 247
$ws.Range("F22").Value = 2082
$ws.Range("F23").Value = 701
$ws.Range("F25").Value = 609
$ws.Range("F26").Value = 374
$ws.Range("F27").Value = 637
$ws.Range("F28").Value = 211
$ws.Range("F29").Value = 222
$ws.Range("F33").Value = 334
$ws.Range("F34").Value = 55
$ws.Range("F36").Value = 226
$ws.Range("F37").Value = 76
$ws.Range("F39").Value = 141
$ws.Range("F40").Value = 40
$ws.Range("F50").Value = 84

